$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 01:51"

# Update country rows whose stats (and consequently sort order) changed
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2833147
$ws.Range("C4").Value = 53194
$ws.Range("D4").Value = 1185477
$ws.Range("E4").Value = 1516257
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 615
$ws.Range("H4").Value = 131413

$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 1501353
$ws.Range("C5").Value = 47984
$ws.Range("D5").Value = 916147
$ws.Range("E5").Value = 523216
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1277
$ws.Range("H5").Value = 61990

$ws.Range("A18").Value = "Alemania"
$ws.Range("B18").Value = 196717
$ws.Range("C18").Value = 393
$ws.Range("D18").Value = 180300
$ws.Range("E18").Value = 7353
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 9064

$ws.Range("A23").Value = "Canada"
$ws.Range("B23").Value = 104771
$ws.Range("C23").Value = 500
$ws.Range("D23").Value = 68347
$ws.Range("E23").Value = 27782
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = 8642

$ws.Range("A27").Value = "Suecia"
$ws.Range("B27").Value = 70639
$ws.Range("C27").Value = 250
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 12
$ws.Range("H27").Value = 5411

$ws.Range("A28").Value = "Argentina"
$ws.Range("B28").Value = 69941
$ws.Range("C28").Value = 2744
$ws.Range("D28").Value = 24186
$ws.Range("E28").Value = 44370
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 34
$ws.Range("H28").Value = 1385

$ws.Range("A43").Value = "Panama"
$ws.Range("B43").Value = 35237
$ws.Range("C43").Value = 774
$ws.Range("D43").Value = 16445
$ws.Range("E43").Value = 18125
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 22
$ws.Range("H43").Value = 667

$ws.Range("A44").Value = "Polonia"
$ws.Range("B44").Value = 35146
$ws.Range("C44").Value = 371
$ws.Range("D44").Value = 22209
$ws.Range("E44").Value = 11445
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 15
$ws.Range("H44").Value = 1492

$ws.Range("A51").Value = "Nigeria"
$ws.Range("B51").Value = 27110
$ws.Range("C51").Value = 626
$ws.Range("D51").Value = 10801
$ws.Range("E51").Value = 15693
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 616

$ws.Range("A52").Value = "Israel"
$ws.Range("B52").Value = 27047
$ws.Range("C52").Value = 790
$ws.Range("D52").Value = 17547
$ws.Range("E52").Value = 9176
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 324

$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 26658
$ws.Range("C53").Value = 593
$ws.Range("D53").Value = 15036
$ws.Range("E53").Value = 11163
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 459

$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 12178
$ws.Range("C69").Value = 132
$ws.Range("D69").Value = 7822
$ws.Range("E69").Value = 4003
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 353

$ws.Range("A137").Value = "Uruguay"
$ws.Range("B137").Value = 947
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 828
$ws.Range("E137").Value = 91
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 28

$ws.Range("A141").Value = "Libia"
$ws.Range("B141").Value = 891
$ws.Range("C141").Value = 17
$ws.Range("D141").Value = 224
$ws.Range("E141").Value = 641
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 26

$ws.Range("A144").Value = "Principado de Andorra"
$ws.Range("B144").Value = 855
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 800
$ws.Range("E144").Value = 3
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 52

$ws.Range("A158").Value = "Vietnam"
$ws.Range("B158").Value = 355
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 340
$ws.Range("E158").Value = 15
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 0

$ws.Range("A184").Value = "Liechtenstein"
$ws.Range("B184").Value = 83
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 81
$ws.Range("E184").Value = 1
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 1

$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 19
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 19
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

$ws.Range("A204").Value = "Laos"
$ws.Range("B204").Value = 19
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 19
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

$ws.Range("A205").Value = "Fiyi"
$ws.Range("B205").Value = 18
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 18
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Dominica"
$ws.Range("B206").Value = 18
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 18
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
